$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.654738187789917
$ws.Range("B1").Value = 2.300026893615723
$ws.Range("C1").Value = 4.429159641265869
$ws.Range("D1").Value = 4.430465698242188
$ws.Range("E1").Value = 1.63140606880188
